$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H11").Value = 1131
$ws.Range("I11").Value = 1131
$ws.Range("K11").Value = 1131
$ws.Range("M11").Value = -991
$ws.Range("H19").Value = 1801.4333
$ws.Range("I19").Value = 1417.3125
$ws.Range("J19").Value = 2240.4285
$ws.Range("K19").Value = 1417.3125
$ws.Range("L19").Value = 2240.4285
$ws.Range("M19").Value = -1242.3125
$ws.Range("N19").Value = -2590.4285
$ws.Range("H64").Value = 6978.5293
$ws.Range("I64").Value = 3549.5557
$ws.Range("J64").Value = 10836.125
$ws.Range("K64").Value = 3549.5557
$ws.Range("L64").Value = 10836.125
$ws.Range("M64").Value = -3301.5557
$ws.Range("N64").Value = -11332.125
$ws.Range("H67").Value = 6978.5293
$ws.Range("I67").Value = 3549.5557
$ws.Range("J67").Value = 10836.125
$ws.Range("K67").Value = 3549.5557
$ws.Range("L67").Value = 10836.125
$ws.Range("M67").Value = -2691.5557
$ws.Range("N67").Value = -12552.125
$ws.Range("H82").Value = 1949.25
$ws.Range("I82").Value = 1949.25
$ws.Range("K82").Value = 5847.75
$ws.Range("M82").Value = -5441.75
$ws.Range("H85").Value = 1949.25
$ws.Range("I85").Value = 1949.25
$ws.Range("K85").Value = 5847.75
$ws.Range("M85").Value = -4443.75
$ws.Range("H86").Value = 22187.75
$ws.Range("I86").Value = 22187.75
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 22187.75
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -21064.75
$ws.Range("H89").Value = 22187.75
$ws.Range("I89").Value = 22187.75
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 110938.75
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -105322.75
$ws.Range("H97").Value = 1893.2222
$ws.Range("J97").Value = 1893.2222
$ws.Range("L97").Value = 5679.6666
$ws.Range("N97").Value = -6671.6666
$ws.Range("H106").Value = 3814
$ws.Range("I106").Value = 3144.8
$ws.Range("J106").Value = 4929.3335
$ws.Range("K106").Value = 3144.8
$ws.Range("L106").Value = 4929.3335
$ws.Range("M106").Value = -2513.8
$ws.Range("N106").Value = -6191.3335
$ws.Range("H107").Value = 477.7143
$ws.Range("I107").Value = 506.9
$ws.Range("J107").Value = 404.75
$ws.Range("K107").Value = 506.9
$ws.Range("L107").Value = 404.75
$ws.Range("M107").Value = 1413.1
$ws.Range("N107").Value = -4244.75
$ws.Range("H129").Value = 4630764.5
$ws.Range("I129").Value = 745.5714
$ws.Range("J129").Value = 20835832
$ws.Range("K129").Value = 2236.7142
$ws.Range("L129").Value = 62507496
$ws.Range("M129").Value = 2763.2858
$ws.Range("N129").Value = -62517496
$ws.Range("H141").Value = 1740
$ws.Range("I141").Value = 1740
$ws.Range("K141").Value = 5220
$ws.Range("M141").Value = -40
$ws.Range("N86").ClearContents()
$ws.Range("N89").ClearContents()

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 4498.7144
$ws.Range("I61").Value = 4498.7144
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4498.7144
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4286.7144
$ws.Range("H74").Value = 5423.65
$ws.Range("I74").Value = 4359.6113
$ws.Range("J74").Value = 15000
$ws.Range("K74").Value = 4359.6113
$ws.Range("L74").Value = 15000
$ws.Range("M74").Value = -3485.6113
$ws.Range("N74").Value = -16748
$ws.Range("H77").Value = 5423.65
$ws.Range("I77").Value = 4359.6113
$ws.Range("J77").Value = 15000
$ws.Range("K77").Value = 21798.0565
$ws.Range("L77").Value = 75000
$ws.Range("M77").Value = -17430.0565
$ws.Range("N77").Value = -83736
$ws.Range("H132").Value = 1806.0416
$ws.Range("I132").Value = 1886.579
$ws.Range("K132").Value = 5659.737
$ws.Range("M132").Value = -3129.737
$ws.Range("H136").Value = 4498.7144
$ws.Range("I136").Value = 4498.7144
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 13496.1432
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -10946.1432
$ws.Range("N61").ClearContents()
$ws.Range("N136").ClearContents()

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H134").Value = 4332.2583
$ws.Range("I134").Value = 3493.8147
$ws.Range("K134").Value = 10481.4441
$ws.Range("M134").Value = -7946.444100000001

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H69").Value = 20007.6
$ws.Range("I69").Value = 13352.667
$ws.Range("J69").Value = 29990
$ws.Range("K69").Value = 13352.667
$ws.Range("L69").Value = 29990
$ws.Range("M69").Value = -12603.667
$ws.Range("N69").Value = -31488
$ws.Range("H72").Value = 20007.6
$ws.Range("I72").Value = 13352.667
$ws.Range("J72").Value = 29990
$ws.Range("K72").Value = 40058.001
$ws.Range("L72").Value = 89970
$ws.Range("M72").Value = -36314.001
$ws.Range("N72").Value = -97458
$ws.Range("H105").Value = 994.8
$ws.Range("I105").Value = 906.125
$ws.Range("K105").Value = 906.125
$ws.Range("M105").Value = 840.875
$ws.Range("H107").Value = 854.15
$ws.Range("I107").Value = 392.35715
$ws.Range("J107").Value = 1931.6666
$ws.Range("K107").Value = 392.35715
$ws.Range("L107").Value = 1931.6666
$ws.Range("M107").Value = 1527.64285
$ws.Range("N107").Value = -5771.6666
$ws.Range("H132").Value = 3717.087
$ws.Range("I132").Value = 2365.4707
$ws.Range("J132").Value = 7546.6665
$ws.Range("K132").Value = 7096.4121
$ws.Range("L132").Value = 22639.9995
$ws.Range("M132").Value = -4566.4121
$ws.Range("N132").Value = -27699.9995
$ws.Range("H134").Value = 3488.1738
$ws.Range("I134").Value = 2408.0625
$ws.Range("J134").Value = 5957
$ws.Range("K134").Value = 7224.1875
$ws.Range("L134").Value = 17871
$ws.Range("M134").Value = -4689.1875
$ws.Range("N134").Value = -22941

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("H129").Value = 694941.9399999999
$ws.Range("J129").Value = 1010049.8
$ws.Range("L129").Value = 3030149.4
$ws.Range("N129").Value = -3040149.4
$ws.Range("H131").Value = 29415114
$ws.Range("I131").Value = 125002100
$ws.Range("J131").Value = 3733.3076
$ws.Range("K131").Value = 375006300
$ws.Range("L131").Value = 11199.9228
$ws.Range("M131").Value = -375001260
$ws.Range("N131").Value = -21279.9228
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 268813.88
$ws.Range("I113").Value = 364655.28
$ws.Range("K113").Value = 364655.28
$ws.Range("M113").Value = -362485.28
$ws.Range("H122").Value = 1761.8695
$ws.Range("I122").Value = 1779.3889
$ws.Range("J122").Value = 1698.8
$ws.Range("K122").Value = 5338.1667
$ws.Range("L122").Value = 5096.4
$ws.Range("M122").Value = -2888.1667
$ws.Range("N122").Value = -9996.4

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H2").Value = 9000
$ws.Range("J2").Value = 9000
$ws.Range("L2").Value = 9000
$ws.Range("N2").Value = -9224
$ws.Range("H10").Value = 3600
$ws.Range("J10").Value = 3600
$ws.Range("L10").Value = 3600
$ws.Range("N10").Value = -3880
$ws.Range("H93").Value = 13436.823
$ws.Range("I93").Value = 1948.5385
$ws.Range("J93").Value = 50773.75
$ws.Range("K93").Value = 1948.5385
$ws.Range("L93").Value = 50773.75
$ws.Range("M93").Value = -700.5385000000001
$ws.Range("N93").Value = -53269.75
$ws.Range("H136").Value = 117650920
$ws.Range("I136").Value = 76926390
$ws.Range("K136").Value = 230779170
$ws.Range("M136").Value = -230776620

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H3").Value = 10431
$ws.Range("J3").Value = 647
$ws.Range("L3").Value = 647
$ws.Range("N3").Value = -875
$ws.Range("H54").Value = 61256.668
$ws.Range("J54").Value = 61256.668
$ws.Range("L54").Value = 61256.668
$ws.Range("N54").Value = -62296.668
$ws.Range("H81").Value = 4914
$ws.Range("I81").Value = 3028
$ws.Range("J81").Value = 7554.4
$ws.Range("K81").Value = 6056
$ws.Range("L81").Value = 15108.8
$ws.Range("M81").Value = -4995
$ws.Range("N81").Value = -17230.8
$ws.Range("H84").Value = 4914
$ws.Range("I84").Value = 3028
$ws.Range("J84").Value = 7554.4
$ws.Range("K84").Value = 30280
$ws.Range("L84").Value = 75544
$ws.Range("M84").Value = -24976
$ws.Range("N84").Value = -86152
$ws.Range("H113").Value = 1100.3334
$ws.Range("I113").Value = 865.56525
$ws.Range("J113").Value = 2450.25
$ws.Range("K113").Value = 2596.69575
$ws.Range("L113").Value = 7350.75
$ws.Range("M113").Value = -426.6957499999999
$ws.Range("N113").Value = -11690.75
$ws.Range("H132").Value = 4269.875
$ws.Range("I132").Value = 3732.4412
$ws.Range("K132").Value = 11197.3236
$ws.Range("M132").Value = -8667.3236
$ws.Range("H136").Value = 3867.6553
$ws.Range("I136").Value = 2079.1904
$ws.Range("J136").Value = 8562.375
$ws.Range("K136").Value = 6237.5712
$ws.Range("L136").Value = 25687.125
$ws.Range("M136").Value = -3867.5712
$ws.Range("N136").Value = -30787.125
